# Auto-generated PowerShell COM-interop script
# Adds a new '2022-Q3' sheet with fund-holding data, right after the '总计' (total) sheet,
# and inserts a corresponding new first-data-row into the '总计' summary sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---- 1. Create the new '2022-Q3' worksheet right after '总计' ----
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Pre-format the text-like numeric columns (B:G) as Text so codes/figures
# such as "014772" or "5.05" are stored as strings, matching the source data,
# instead of being auto-coerced into numbers.
$q3.Range("B2:G15").NumberFormat = "@"

# Data rows
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "014772"
$q3.Cells.Item(2,3).Value = "中泰红利价值一年持有混合"
$q3.Cells.Item(2,4).Value = "5.05"
$q3.Cells.Item(2,5).Value = "93.55"
$q3.Cells.Item(2,6).Value = "4.66"
$q3.Cells.Item(2,7).Value = "0.2353"
$q3.Cells.Item(2,8).Value = 10

$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "862001"
$q3.Cells.Item(3,3).Value = "光大阳光香港精选混合（QDII）A 人民币"
$q3.Cells.Item(3,4).Value = "3.15"
$q3.Cells.Item(3,5).Value = "90.62"
$q3.Cells.Item(3,6).Value = "4.47"
$q3.Cells.Item(3,7).Value = "0.1408"
$q3.Cells.Item(3,8).Value = 7

$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "862011"
$q3.Cells.Item(4,3).Value = "光大阳光香港精选混合（QDII）A 美元"
$q3.Cells.Item(4,4).Value = "3.15"
$q3.Cells.Item(4,5).Value = "90.62"
$q3.Cells.Item(4,6).Value = "4.47"
$q3.Cells.Item(4,7).Value = "0.1408"
$q3.Cells.Item(4,8).Value = 7

$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "862012"
$q3.Cells.Item(5,3).Value = "光大阳光香港精选混合（QDII）C 人民币"
$q3.Cells.Item(5,4).Value = "3.15"
$q3.Cells.Item(5,5).Value = "90.62"
$q3.Cells.Item(5,6).Value = "4.47"
$q3.Cells.Item(5,7).Value = "0.1408"
$q3.Cells.Item(5,8).Value = 7

$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "013200"
$q3.Cells.Item(6,3).Value = "南方均衡优选一年持有期混合A"
$q3.Cells.Item(6,4).Value = "7.27"
$q3.Cells.Item(6,5).Value = "40.60"
$q3.Cells.Item(6,6).Value = "1.71"
$q3.Cells.Item(6,7).Value = "0.1243"
$q3.Cells.Item(6,8).Value = 5

$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "860027"
$q3.Cells.Item(7,3).Value = "光大阳光价值30个月持有期混合B"
$q3.Cells.Item(7,4).Value = "2.11"
$q3.Cells.Item(7,5).Value = "91.42"
$q3.Cells.Item(7,6).Value = "3.92"
$q3.Cells.Item(7,7).Value = "0.0827"
$q3.Cells.Item(7,8).Value = 6

$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).Value = "513690"
$q3.Cells.Item(8,3).Value = "博时恒生港股通高股息率ETF"
$q3.Cells.Item(8,4).Value = "3.05"
$q3.Cells.Item(8,5).Value = "97.26"
$q3.Cells.Item(8,6).Value = "2.43"
$q3.Cells.Item(8,7).Value = "0.0741"
$q3.Cells.Item(8,8).Value = 8

$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).Value = "010230"
$q3.Cells.Item(9,3).Value = "南方宝昌混合A"
$q3.Cells.Item(9,4).Value = "16.40"
$q3.Cells.Item(9,5).Value = "22.04"
$q3.Cells.Item(9,6).Value = "0.42"
$q3.Cells.Item(9,7).Value = "0.0689"
$q3.Cells.Item(9,8).Value = 10

$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).Value = "005741"
$q3.Cells.Item(10,3).Value = "南方君信灵活配置混合A"
$q3.Cells.Item(10,4).Value = "2.80"
$q3.Cells.Item(10,5).Value = "71.98"
$q3.Cells.Item(10,6).Value = "1.59"
$q3.Cells.Item(10,7).Value = "0.0445"
$q3.Cells.Item(10,8).Value = 10

$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).Value = "860007"
$q3.Cells.Item(11,3).Value = "光大阳光价值30个月持有期混合A"
$q3.Cells.Item(11,4).Value = "0.59"
$q3.Cells.Item(11,5).Value = "91.42"
$q3.Cells.Item(11,6).Value = "3.92"
$q3.Cells.Item(11,7).Value = "0.0231"
$q3.Cells.Item(11,8).Value = 6

$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).Value = "159726"
$q3.Cells.Item(12,3).Value = "华夏恒生中国内地企业高股息率ETF"
$q3.Cells.Item(12,4).Value = "0.84"
$q3.Cells.Item(12,5).Value = "96.48"
$q3.Cells.Item(12,6).Value = "2.68"
$q3.Cells.Item(12,7).Value = "0.0225"
$q3.Cells.Item(12,8).Value = 6

$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).Value = "013201"
$q3.Cells.Item(13,3).Value = "南方均衡优选一年持有期混合C"
$q3.Cells.Item(13,4).Value = "0.73"
$q3.Cells.Item(13,5).Value = "40.60"
$q3.Cells.Item(13,6).Value = "1.71"
$q3.Cells.Item(13,7).Value = "0.0125"
$q3.Cells.Item(13,8).Value = 5

$q3.Cells.Item(14,1).Value = 12
$q3.Cells.Item(14,2).Value = "010150"
$q3.Cells.Item(14,3).Value = "南方君信灵活配置混合C"
$q3.Cells.Item(14,4).Value = "0.20"
$q3.Cells.Item(14,5).Value = "71.98"
$q3.Cells.Item(14,6).Value = "1.59"
$q3.Cells.Item(14,7).Value = "0.0032"
$q3.Cells.Item(14,8).Value = 10

$q3.Cells.Item(15,1).Value = 13
$q3.Cells.Item(15,2).Value = "010231"
$q3.Cells.Item(15,3).Value = "南方宝昌混合C"
$q3.Cells.Item(15,4).Value = "0.53"
$q3.Cells.Item(15,5).Value = "22.04"
$q3.Cells.Item(15,6).Value = "0.42"
$q3.Cells.Item(15,7).Value = "0.0022"
$q3.Cells.Item(15,8).Value = 10

# Bold + centered + thin-bordered style for the header row and the row-index column (A),
# matching the look of the '总计' sheet and the other quarterly sheets.
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$idxRange = $q3.Range("A2:A15")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1
$idxRange.Borders.Weight = 2

# ---- 2. Update the '总计' (total) summary sheet: shift existing rows down
#         by one and insert the new 2022-Q3 figures at the top ----

# Row 9 doesn't exist yet; give A9 the same bold/border/center look as A2:A8
# before writing into it (writing into an untouched cell starts from the plain
# default style, unlike the existing rows which already carry that formatting).
$a9 = $totalSheet.Cells.Item(9,1)
$a9.Font.Bold = $true
$a9.HorizontalAlignment = -4108
$a9.VerticalAlignment = -4160
$a9.Borders.LineStyle = 1
$a9.Borders.Weight = 2

# Rewrite rows 9 down to 2 (bottom-up) with their new target contents so the
# existing quarters shift down by one row and 2022-Q3 lands in row 2.
$totalSheet.Cells.Item(9,1).Value = 7
$totalSheet.Cells.Item(9,2).Value = "2020-Q4"
$totalSheet.Cells.Item(9,3).Value = 6
$totalSheet.Cells.Item(9,4).Value = 4.12

$totalSheet.Cells.Item(8,1).Value = 6
$totalSheet.Cells.Item(8,2).Value = "2021-Q1"
$totalSheet.Cells.Item(8,3).Value = 4
$totalSheet.Cells.Item(8,4).Value = 0.16

$totalSheet.Cells.Item(7,1).Value = 5
$totalSheet.Cells.Item(7,2).Value = "2021-Q2"
$totalSheet.Cells.Item(7,3).Value = 2
$totalSheet.Cells.Item(7,4).Value = 0.04

$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(6,2).Value = "2021-Q3"
$totalSheet.Cells.Item(6,3).Value = 12
$totalSheet.Cells.Item(6,4).Value = 1.11

$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(5,2).Value = "2021-Q4"
$totalSheet.Cells.Item(5,3).Value = 5
$totalSheet.Cells.Item(5,4).Value = 0.56

$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2022-Q1"
$totalSheet.Cells.Item(4,3).Value = 17
$totalSheet.Cells.Item(4,4).Value = 2.75

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q2"
$totalSheet.Cells.Item(3,3).Value = 18
$totalSheet.Cells.Item(3,4).Value = 2.19

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q3"
$totalSheet.Cells.Item(2,3).Value = 14
$totalSheet.Cells.Item(2,4).Value = 1.12

